# Updates the quick-links list on Sheet1:
#  - "NBC Store" is replaced by "Contact Us"
#  - "Accessibility" is replaced by a new entry "Tickets and NBC Studio Tour"
#  - "Contact Us" (previously below) becomes "Accessibility"
# Net effect: a new row "Tickets and NBC Studio Tour" is inserted after the
# "Parental Guidelines and TV Ratings" row, the rest shifting down, and the
# final selection is left on cell A19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Contact Us"
$ws.Range("A9").Value = "Tickets and NBC Studio Tour"
$ws.Range("A10").Value = "Accessibility"

$ws.Range("A19").Select()
